# Apply the "increased title length" content edits to initial_epic_data.xlsx
# - Row 23, col C ("Programs"): "Watershed Planning" -> "Watershed Management"
# - Rows 32-35, col B ("Program areas"): "Drought Monitoring" -> "Drought Monitoring, Response and Recovery"
# - Rows 32-35, col C ("Programs"): leading-space variants -> trimmed text
#     Row 32: " Drought Monitoring Program" -> "Drought Monitoring Program"
#     Row 33: " WRM Drought Response" -> "WRM Drought Response"
#     Row 34: " Agriculture Drought Response" -> "Agriculture Drought Response"
#     Row 35: " Social Protection Drought Response" -> "Social Protection Drought Response"
# - Rows 36-38, col B ("Program areas"): "Flood Monitoring" -> "Flood Monitoring, Response and Recovery"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(23, 3).Value = "Watershed Management"

$ws.Cells.Item(32, 2).Value = "Drought Monitoring, Response and Recovery"
$ws.Cells.Item(33, 2).Value = "Drought Monitoring, Response and Recovery"
$ws.Cells.Item(34, 2).Value = "Drought Monitoring, Response and Recovery"
$ws.Cells.Item(35, 2).Value = "Drought Monitoring, Response and Recovery"

$ws.Cells.Item(32, 3).Value = "Drought Monitoring Program"
$ws.Cells.Item(33, 3).Value = "WRM Drought Response"
$ws.Cells.Item(34, 3).Value = "Agriculture Drought Response"
$ws.Cells.Item(35, 3).Value = "Social Protection Drought Response"

$ws.Cells.Item(36, 2).Value = "Flood Monitoring, Response and Recovery"
$ws.Cells.Item(37, 2).Value = "Flood Monitoring, Response and Recovery"
$ws.Cells.Item(38, 2).Value = "Flood Monitoring, Response and Recovery"

# Update the saved view/selection state to match the author's last interaction:
# top-left visible cell A16, active selection on B40.
$ws.Range("A16").Select()
$ws.Range("B40").Select()
